$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")
$cols = @("D","E","F","G","H","I","J","K","L","M")

$row8 = @("9 ماهه منتهی به 1399/09", "12 ماهه منتهی به 1399/12", "3 ماهه منتهی به 1400/03", "6 ماهه منتهی به 1400/06", "9 ماهه منتهی به 1400/09", "12 ماهه منتهی به 1400/12", "3 ماهه منتهی به 1401/03", "6 ماهه منتهی به 1401/06", "9 ماهه منتهی به 1401/09", "12 ماهه منتهی به 1401/12")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "8").Value = $row8[$i]
}

$row9 = @("1401-01-15 (3)", "1401-03-11 (10)", "1401-04-30 (3)", "1401-08-18 (4)", "1401-10-29 (3)", "1402-02-27 (7)", "1401-04-30", "1401-08-18 (2)", "1401-10-29", "1402-02-27")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "9").Value = $row9[$i]
}

$row11 = @(8847, 10796, 4105, 8826, 13413, 16287, 4119, 8810, 19985, 27827)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "11").Value = $row11[$i]
}

$row12 = @(-4689, -5309, -1999, -4151, -6753, -8064, -2051, -4086, -11288, -16737)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "12").Value = $row12[$i]
}

$row13 = @(4158, 5487, 2106, 4675, 6661, 8223, 2068, 4724, 8698, 11090)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "13").Value = $row13[$i]
}

$row14 = @(-286, -516, -165, -333, -453, -674, -356, -364, -539, -606)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "14").Value = $row14[$i]
}

$row15 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "15").Value = $row15[$i]
}

$row16 = @(75, 69, 5, 73, 113, 153, 0, 19, 75, 81)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "16").Value = $row16[$i]
}

$row17 = @(3948, 5040, 1946, 4415, 6321, 7702, 1713, 4378, 8234, 10564)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "17").Value = $row17[$i]
}

$row18 = @(-404, -688, -117, -384, -616, -881, -88, -538, -736, -956)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "18").Value = $row18[$i]
}

$row19 = @(2, -310, 0, -125, "-", -127, 29, 27, 28, -7)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "19").Value = $row19[$i]
}

$row20 = @(3546, 4043, 1828, 3907, 5705, 6694, 1653, 3867, 7526, 9601)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "20").Value = $row20[$i]
}

$row21 = @(-806, -487, -416, -875, -1188, -943, -372, -870, -1693, -1422)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "21").Value = $row21[$i]
}

$row22 = @(2740, 3556, 1413, 3031, 4517, 5751, 1281, 2997, 5832, 8179)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "22").Value = $row22[$i]
}

$row23 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "23").Value = $row23[$i]
}

$row24 = @(2740, 3556, 1413, 3031, 4517, 5751, 1281, 2997, 5832, 8179)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "24").Value = $row24[$i]
}

$row25 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "25").Value = $row25[$i]
}

$row26 = @(1159, 1128, 1950, 1836, 1749, 1724, 1548, 1507, 1431, 2550)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "26").Value = $row26[$i]
}

$row27 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "27").Value = $row27[$i]
}
